# ============================================================
# Update Daily Report: 2026-02-03
# Adds the 2026-02-02 trading-day rows to Daily_Data, and
# rolls Monthly_Stats forward to include a 2026-02 month bucket
# (summary row + per-depository detail rows), pushing the
# existing 2026-01 rows down.
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Daily_Data: append the new day's 22 rows (rows 442-463)
# ------------------------------------------------------------------
$dd = $wb.Worksheets.Item("Daily_Data")

$newDate = 46055

$dailyRows = @(
    @("ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @("ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("BRINK'S, INC. Registered", 76497.842, 0, 0, 0, 0, 76497.842),
    @("BRINK'S, INC. Eligible", 42030.257, 0, 0, 0, 0, 42030.257),
    @("CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @("CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @("DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @("DELAWARE DEPOSITORY Eligible", 18459.584, 0, 0, 0, 0, 18459.584),
    @("HSBC BANK, USA Registered", 1394.758, 0, 0, 0, 0, 1394.758),
    @("HSBC BANK, USA Eligible", 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @("JP MORGAN CHASE BANK NA Registered", 114985.579, 0, 0, 0, 0, 114985.579),
    @("JP MORGAN CHASE BANK NA Eligible", 75484.511, 0, 0, 0, 0, 75484.511),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 63745.991, 0, 0, 0, 0, 63745.991),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 132077.206, 0, 0, 0, 0, 132077.206),
    @("MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @("MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 50220.42, 0, 0, 0, 0, 50220.42),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 1271.373, 0, 0, 0, 0, 1271.373),
    @("STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @("STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

$startRow = 442
for ($i = 0; $i -lt $dailyRows.Length; $i++) {
    $row = $dailyRows[$i]
    $r = $startRow + $i

    $dateCell = $dd.Cells.Item($r, 1)
    $dateCell.Value = $newDate
    $dateCell.NumberFormat = $dd.Cells.Item($r - 1, 1).NumberFormat

    $dd.Cells.Item($r, 2).Value = $row[0]
    $dd.Cells.Item($r, 3).Value = $row[1]
    $dd.Cells.Item($r, 4).Value = $row[2]
    $dd.Cells.Item($r, 5).Value = $row[3]
    $dd.Cells.Item($r, 6).Value = $row[4]
    $dd.Cells.Item($r, 7).Value = $row[5]
    $dd.Cells.Item($r, 8).Value = $row[6]
}

# ------------------------------------------------------------------
# 2) Monthly_Stats: roll the month forward
# ------------------------------------------------------------------
$ms = $wb.Worksheets.Item("Monthly_Stats")

$shiftDown = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown
$fromBelow = [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromRightOrBelow

# 2a) top summary block: row 2 becomes 2026-02 (values unchanged),
#     a new row 3 is inserted carrying the old 2026-01 summary values.
$oldB2 = $ms.Cells.Item(2, 2).Value2
$oldC2 = $ms.Cells.Item(2, 3).Value2
$oldD2 = $ms.Cells.Item(2, 4).Value2

$ms.Cells.Item(2, 1).Value = "2026-02"

$ms.Range("A3:A3").EntireRow.Insert($shiftDown, $fromBelow)
$ms.Cells.Item(3, 1).Value = "2026-01"
$ms.Cells.Item(3, 2).Value = $oldB2
$ms.Cells.Item(3, 3).Value = $oldC2
$ms.Cells.Item(3, 4).Value = $oldD2

# 2b) detail block: insert 22 fresh rows right after the header
#     (old row 6 -> 7) for the new 2026-02 per-depository entries;
#     the pre-existing 2026-01 detail rows are pushed down intact.
$ms.Range("A8:A29").EntireRow.Insert($shiftDown, $fromBelow)

$febRows = @(
    @("ASAHI DEPOSITORY LLC Eligible", 0, 0, 0),
    @("ASAHI DEPOSITORY LLC Registered", 0, 0, 0),
    @("BRINK'S, INC. Eligible", 0, 0, 42030.257),
    @("BRINK'S, INC. Registered", 0, 0, 76497.842),
    @("CNT DEPOSITORY, INC. Eligible", 0, 0, 0),
    @("CNT DEPOSITORY, INC. Registered", 0, 0, 1246.06),
    @("DELAWARE DEPOSITORY Eligible", 0, 0, 18459.584),
    @("DELAWARE DEPOSITORY Registered", 0, 0, 1633.941),
    @("HSBC BANK, USA Eligible", 0, 0, 9281.979),
    @("HSBC BANK, USA Registered", 0, 0, 1394.758),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 0, 0, 2395.448),
    @("JP MORGAN CHASE BANK NA Eligible", 0, 0, 75484.511),
    @("JP MORGAN CHASE BANK NA Registered", 0, 0, 114985.579),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 0, 0, 132077.206),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 0, 0, 63745.991),
    @("MALCA-AMIT USA, LLC Eligible", 0, 0, 0),
    @("MALCA-AMIT USA, LLC Registered", 0, 0, 395.145),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 0, 0, 1271.373),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 0, 0, 50220.42),
    @("STONEX PRECIOUS METALS LLC Eligible", 0, 0, 16.075),
    @("STONEX PRECIOUS METALS LLC Registered", 0, 0, 14122.765)
)

$startRow2 = 8
for ($i = 0; $i -lt $febRows.Length; $i++) {
    $row = $febRows[$i]
    $r = $startRow2 + $i

    $ms.Cells.Item($r, 1).Value = "2026-02"
    $ms.Cells.Item($r, 2).Value = $row[0]
    $ms.Cells.Item($r, 3).Value = $row[1]
    $ms.Cells.Item($r, 4).Value = $row[2]
    $ms.Cells.Item($r, 5).Value = $row[3]
}
